$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 94
$ws.Cells.Item(94, 8).Value = 1835
$ws.Cells.Item(94, 9).Value = 1835
$ws.Cells.Item(94, 11).Value = 1835
$ws.Cells.Item(94, 13).Value = -1384
# Row 113
$ws.Cells.Item(113, 8).Value = 2309.1924
$ws.Cells.Item(113, 9).Value = 2183.6365
$ws.Cells.Item(113, 10).Value = 2999.75
$ws.Cells.Item(113, 11).Value = 2183.6365
$ws.Cells.Item(113, 12).Value = 2999.75
$ws.Cells.Item(113, 13).Value = 1070.3635
$ws.Cells.Item(113, 14).Value = -9507.75
# Row 132
$ws.Cells.Item(132, 8).Value = 4934.4443
$ws.Cells.Item(132, 9).Value = 5202.7827
$ws.Cells.Item(132, 11).Value = 15608.3481
$ws.Cells.Item(132, 13).Value = -13078.3481

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 5221
$ws.Cells.Item(2, 9).Value = 2046
$ws.Cells.Item(2, 10).Value = 7337.6665
$ws.Cells.Item(2, 11).Value = 2046
$ws.Cells.Item(2, 12).Value = 7337.6665
$ws.Cells.Item(2, 13).Value = -1933
$ws.Cells.Item(2, 14).Value = -7563.6665
# Row 32
$ws.Cells.Item(32, 8).Value = 417328.9
$ws.Cells.Item(32, 9).Value = 2662.1865
$ws.Cells.Item(32, 10).Value = 2863862.5
$ws.Cells.Item(32, 11).Value = 2662.1865
$ws.Cells.Item(32, 12).Value = 2863862.5
$ws.Cells.Item(32, 13).Value = -2375.1865
$ws.Cells.Item(32, 14).Value = -2864436.5
# Row 61
$ws.Cells.Item(61, 8).Value = 1308089.1
$ws.Cells.Item(61, 9).Value = 1667465.2
$ws.Cells.Item(61, 10).Value = 1267.2727
$ws.Cells.Item(61, 11).Value = 1667465.2
$ws.Cells.Item(61, 12).Value = 1267.2727
$ws.Cells.Item(61, 13).Value = -1667253.2
$ws.Cells.Item(61, 14).Value = -1691.2727
# Row 116
$ws.Cells.Item(116, 8).Value = 5221
$ws.Cells.Item(116, 9).Value = 2046
$ws.Cells.Item(116, 10).Value = 7337.6665
$ws.Cells.Item(116, 11).Value = 2046
$ws.Cells.Item(116, 12).Value = 7337.6665
$ws.Cells.Item(116, 13).Value = 248
$ws.Cells.Item(116, 14).Value = -11925.6665
# Row 123
$ws.Cells.Item(123, 8).Value = 54980
$ws.Cells.Item(123, 10).Value = 54980
$ws.Cells.Item(123, 12).Value = 54980
$ws.Cells.Item(123, 14).Value = -64780
# Row 132
$ws.Cells.Item(132, 8).Value = 29018.305
$ws.Cells.Item(132, 9).Value = 773.8077
$ws.Cells.Item(132, 10).Value = 102454
$ws.Cells.Item(132, 11).Value = 2321.4231
$ws.Cells.Item(132, 12).Value = 307362
$ws.Cells.Item(132, 13).Value = 208.5769
$ws.Cells.Item(132, 14).Value = -312422
# Row 136
$ws.Cells.Item(136, 8).Value = 1308089.1
$ws.Cells.Item(136, 9).Value = 1667465.2
$ws.Cells.Item(136, 10).Value = 1267.2727
$ws.Cells.Item(136, 11).Value = 5002395.6
$ws.Cells.Item(136, 12).Value = 3801.8181
$ws.Cells.Item(136, 13).Value = -4999845.6
$ws.Cells.Item(136, 14).Value = -8901.8181

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 5221
$ws.Cells.Item(3, 9).Value = 2046
$ws.Cells.Item(3, 10).Value = 7337.6665
$ws.Cells.Item(3, 11).Value = 2046
$ws.Cells.Item(3, 12).Value = 7337.6665
$ws.Cells.Item(3, 13).Value = -1932
$ws.Cells.Item(3, 14).Value = -7565.6665
# Row 134
$ws.Cells.Item(134, 8).Value = 6143.162
$ws.Cells.Item(134, 9).Value = 2042.6897
$ws.Cells.Item(134, 10).Value = 21007.375
$ws.Cells.Item(134, 11).Value = 6128.0691
$ws.Cells.Item(134, 12).Value = 63022.125
$ws.Cells.Item(134, 13).Value = -3593.0691
$ws.Cells.Item(134, 14).Value = -68092.125

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 5103857
$ws.Cells.Item(31, 9).Value = 6758179.5
$ws.Cells.Item(31, 10).Value = 3029.0833
$ws.Cells.Item(31, 11).Value = 6758179.5
$ws.Cells.Item(31, 12).Value = 3029.0833
$ws.Cells.Item(31, 13).Value = -6757884.5
$ws.Cells.Item(31, 14).Value = -3619.0833
# Row 34
$ws.Cells.Item(34, 8).Value = 5103857
$ws.Cells.Item(34, 9).Value = 6758179.5
$ws.Cells.Item(34, 10).Value = 3029.0833
$ws.Cells.Item(34, 11).Value = 6758179.5
$ws.Cells.Item(34, 12).Value = 3029.0833
$ws.Cells.Item(34, 13).Value = -6757977.5
$ws.Cells.Item(34, 14).Value = -3433.0833
# Row 86
$ws.Cells.Item(86, 8).Value = 52618.5
$ws.Cells.Item(86, 9).Value = 18150
$ws.Cells.Item(86, 11).Value = 18150
$ws.Cells.Item(86, 13).Value = -17027
# Row 89
$ws.Cells.Item(89, 8).Value = 52618.5
$ws.Cells.Item(89, 9).Value = 18150
$ws.Cells.Item(89, 11).Value = 90750
$ws.Cells.Item(89, 13).Value = -85134
# Row 132
$ws.Cells.Item(132, 8).Value = 42186.64
$ws.Cells.Item(132, 9).Value = 1459.8572
$ws.Cells.Item(132, 10).Value = 58024.832
$ws.Cells.Item(132, 11).Value = 4379.571599999999
$ws.Cells.Item(132, 12).Value = 174074.496
$ws.Cells.Item(132, 13).Value = -1849.571599999999
$ws.Cells.Item(132, 14).Value = -179134.496

$ws = $wb.Worksheets.Item("CUL")
# Row 93
$ws.Cells.Item(93, 8).Value = 4300
$ws.Cells.Item(93, 10).Value = 4300
$ws.Cells.Item(93, 12).Value = 12900
$ws.Cells.Item(93, 14).Value = -16644
# Row 94
$ws.Cells.Item(94, 8).Value = 6131
$ws.Cells.Item(94, 9).Value = 4262.3335
$ws.Cells.Item(94, 10).Value = 6385.8184
$ws.Cells.Item(94, 11).Value = 12787.0005
$ws.Cells.Item(94, 12).Value = 19157.4552
$ws.Cells.Item(94, 13).Value = -12111.0005
$ws.Cells.Item(94, 14).Value = -20509.4552
# Row 95
$ws.Cells.Item(95, 8).Value = 17500
$ws.Cells.Item(95, 9).Value = 30000
$ws.Cells.Item(95, 10).Value = 5000
$ws.Cells.Item(95, 11).Value = 90000
$ws.Cells.Item(95, 12).Value = 15000
$ws.Cells.Item(95, 13).Value = -87941
$ws.Cells.Item(95, 14).Value = -19118
# Row 97
$ws.Cells.Item(97, 8).Value = 398.81818
$ws.Cells.Item(97, 9).Value = 513.1667
$ws.Cells.Item(97, 10).Value = 261.6
$ws.Cells.Item(97, 11).Value = 1539.5001
$ws.Cells.Item(97, 12).Value = 784.8000000000001
$ws.Cells.Item(97, 13).Value = -1043.5001
$ws.Cells.Item(97, 14).Value = -1776.8
# Row 101
$ws.Cells.Item(101, 8).Value = 5357.143
$ws.Cells.Item(101, 10).Value = 5357.143
$ws.Cells.Item(101, 12).Value = 16071.429
$ws.Cells.Item(101, 14).Value = -20939.429
# Row 131
$ws.Cells.Item(131, 8).Value = 23810792
$ws.Cells.Item(131, 9).Value = 831.5
$ws.Cells.Item(131, 10).Value = 55557404
$ws.Cells.Item(131, 11).Value = 2494.5
$ws.Cells.Item(131, 12).Value = 166672212
$ws.Cells.Item(131, 13).Value = 2545.5
$ws.Cells.Item(131, 14).Value = -166682292

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 2136
$ws.Cells.Item(40, 9).Value = 2128.4285
$ws.Cells.Item(40, 11).Value = 2128.4285
$ws.Cells.Item(40, 13).Value = -1992.4285
# Row 132
$ws.Cells.Item(132, 8).Value = 78618.71000000001
$ws.Cells.Item(132, 9).Value = 84281.69500000001
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 252845.085
$ws.Cells.Item(132, 12).Value = 15000
$ws.Cells.Item(132, 13).Value = -250315.085
$ws.Cells.Item(132, 14).Value = -20060
# Row 136
$ws.Cells.Item(136, 8).Value = 5108.077
$ws.Cells.Item(136, 9).Value = 4972.148
$ws.Cells.Item(136, 10).Value = 5413.9165
$ws.Cells.Item(136, 11).Value = 14916.444
$ws.Cells.Item(136, 12).Value = 16241.7495
$ws.Cells.Item(136, 13).Value = -12366.444
$ws.Cells.Item(136, 14).Value = -21341.7495

$ws = $wb.Worksheets.Item("WVR")
# Row 123
$ws.Cells.Item(123, 8).Value = 37607.43
$ws.Cells.Item(123, 10).Value = 37607.43
$ws.Cells.Item(123, 12).Value = 37607.43
$ws.Cells.Item(123, 14).Value = -47407.43
